$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.926.15"
$ws.Range("E2").Value = "  +6.94%  "

$ws.Range("D3").Value = "3.016.64"
$ws.Range("E3").Value = "  +4.24%  "

$ws.Range("E4").Value = "  +0.16%  "

$ws.Range("D5").Value = "'585.06"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.10%  "

$ws.Range("D6").Value = "'155.94"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +8.88%  "

$ws.Range("D7").Value = "'0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.15%  "

$ws.Range("D8").Value = "3.011.47"
$ws.Range("E8").Value = "  +4.15%  "

$ws.Range("E9").Value = "  +3.18%  "

$ws.Range("D10").Value = "'6.97"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.35%  "

$ws.Range("E11").Value = "  +7.17%  "

$ws.Range("E12").Value = "  +5.33%  "

$ws.Range("D13").Value = "'0.0000255"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +10.25%  "

$ws.Range("E14").Value = "  +8.58%  "

$ws.Range("E15").Value = "  +0.46%  "

$ws.Range("D16").Value = "65.987.95"
$ws.Range("E16").Value = "  +7.05%  "

$ws.Range("D17").Value = "3.518.48"
$ws.Range("E17").Value = "  +4.27%  "

$ws.Range("E18").Value = "  +6.32%  "

$ws.Range("D19").Value = "3.025.63"
$ws.Range("E19").Value = "  +4.45%  "

$ws.Range("D20").Value = "'465.02"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +7.75%  "

$ws.Range("D21").Value = "'13.83"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +6.14%  "

$ws.Range("D22").Value = "'0.683"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +4.57%  "

$ws.Range("E23").Value = "  +8.11%  "

$ws.Range("D24").Value = "'82.10"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.75%  "

$ws.Range("D25").Value = "'12.49"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +5.25%  "

$ws.Range("D26").Value = "'2.24"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +11.48%  "

$ws.Range("D27").Value = "'10.72"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +7.85%  "

$ws.Range("D29").Value = "'7.96"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +13.31%  "

$ws.Range("E30").Value = "  +17.72%  "

$ws.Range("E31").Value = "  -0.29%  "

$ws.Range("E32").Value = "  +4.91%  "

$ws.Range("E33").Value = "  +5.33%  "

$ws.Range("D34").Value = "'27.06"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +6.11%  "

$ws.Range("D35").Value = "'1.00"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.07%  "

$ws.Range("D36").Value = "'0.999"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.25%  "

$ws.Range("D37").Value = "'5.82"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +8.52%  "

$ws.Range("E38").Value = "  +12.88%  "

$ws.Range("D39").Value = "'3.03"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +8.06%  "

$ws.Range("D40").Value = "'49.24"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.92%  "

$ws.Range("D41").Value = "'44.61"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +12.35%  "

$ws.Range("E42").Value = "  +8.21%  "

$ws.Range("D43").Value = "'0.302"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +13.22%  "

$ws.Range("D44").Value = "'8.48"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.61%  "

$ws.Range("D45").Value = "'395.88"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +15.29%  "

$ws.Range("D46").Value = "2.797.41"
$ws.Range("E46").Value = "  +4.01%  "

$ws.Range("D47").Value = "'0.0354"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +5.66%  "

$ws.Range("D48").Value = "'134.09"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.11%  "

$ws.Range("D50").Value = "'23.66"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +9.91%  "

$ws.Range("E51").Value = "  +4.08%  "
